{"js": "// Fixed #348 Hyperlinks from sub-template does not work.\n//\n// The second paragraph of the document (\"name = World,\") still carries the\n// \"blank placeholder run\" left over from the template engine (a <w:pPr> with\n// an en-US language mark, followed by an empty run that also carries that\n// language mark) instead of the normal Word \"proofing error\" bookmarks\n// (spellStart/spellEnd/gramStart/gramEnd) that every other \"name = ...\"\n// paragraph in the document already has. This routine locates that\n// paragraph, drops the stray paragraph-mark formatting / empty run, and\n// inserts the four <w:proofErr/> markers immediately before the existing\n// \"name\" run - bringing the paragraph in line with its siblings.\n\nconst body = context.document.body;\n\n// Locate the paragraph holding the literal text \"name = World,\" - this is\n// unique in the document, so it unambiguously identifies our target even if\n// paragraph indices ever shift.\nconst found = body.search(\"name = World,\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length === 0) {\n  throw new Error(\"Could not locate target paragraph 'name = World,'\");\n}\n\nconst targetParagraph = found.items[0].paragraphs.getFirst();\n\n// Pull the paragraph's own OOXML so we can recover its real w:rsid*\n// attributes (they must be preserved - only the inner content changes).\nconst ooxml = targetParagraph.getOoxml();\nawait context.sync();\n\nconst openTagMatch = ooxml.value.match(/<w:p\\b([^>]*)>/);\nlet pAttrs = openTagMatch ? openTagMatch[1] : \"\";\n// The OOXML round-trip mints fresh w14:paraId/w14:textId attributes that\n// were not present on the original paragraph - strip those back out so we\n// only keep the genuine rsid* attributes.\npAttrs = pAttrs\n  .replace(/\\s*w14:paraId=\"[^\"]*\"/, \"\")\n  .replace(/\\s*w14:textId=\"[^\"]*\"/, \"\")\n  .trim();\nconst pOpenTag = pAttrs ? `<w:p ${pAttrs}>` : \"<w:p>\";\n\n// Rebuild the paragraph: four proofErr markers followed by the same run\n// content as before (name / \" \" / \"= \" / \"World\" / \",\") but without the\n// leftover pPr + empty run.\nconst replacementOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  pOpenTag +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  \"<w:r><w:t>name</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">= </w:t></w:r>' +\n  \"<w:r><w:t>World</w:t></w:r>\" +\n  \"<w:r><w:t>,</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nconst wholeParagraphRange = targetParagraph.getRange(\"Whole\");\nwholeParagraphRange.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Fixed #348 Hyperlinks from sub-template does not work.\n#\n# The second paragraph of the document (\"name = World,\") still carries the\n# stray \"blank placeholder run\" left over from the template engine (a\n# <w:pPr> with an en-US language mark, followed by an empty run that also\n# carries that language mark) instead of the normal Word \"proofing error\"\n# bookmarks (spellStart/spellEnd/gramStart/gramEnd) that every other\n# \"name = ...\" paragraph in the document already has. This script locates\n# that paragraph, drops the stray paragraph-mark formatting / empty run,\n# and inserts the four proofErr markers immediately before the existing\n# \"name\" run - bringing the paragraph in line with its siblings.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph holding the literal text \"name = World,\" - this is\n# unique in the document, so it unambiguously identifies our target even if\n# paragraph indices ever shift. Walking Paragraphs directly (rather than\n# collapsing Document.Content with Find.Execute) keeps the resulting Range\n# anchored to the real paragraph extent, which InsertXML later replaces.\n$targetParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains(\"name\") -and $p.Range.Text.Contains(\"World\")) {\n        $targetParagraph = $p\n        break\n    }\n}\nif ($null -eq $targetParagraph) {\n    throw \"Could not locate target paragraph 'name = World,'\"\n}\n\n$target = $targetParagraph.Range\n\n# Pull the paragraph's own OOXML so we can recover its real w:rsid*\n# attributes (they must be preserved - only the inner content changes).\n$xml = $target.WordOpenXML\nif ($xml -match '<w:p\\s([^>]*)>') {\n    $pAttrs = $matches[1]\n} else {\n    $pAttrs = ''\n}\n# The OOXML round-trip mints fresh w14:paraId/w14:textId attributes that\n# were not present on the original paragraph - strip those back out so we\n# only keep the genuine rsid* attributes.\n$pAttrs = $pAttrs -replace 'w14:paraId=\"[^\"]*\"\\s*', ''\n$pAttrs = $pAttrs -replace 'w14:textId=\"[^\"]*\"\\s*', ''\n$pAttrs = $pAttrs.Trim()\nif ($pAttrs -ne '') {\n    $pOpenTag = \"<w:p $pAttrs>\"\n} else {\n    $pOpenTag = '<w:p>'\n}\n\n# Rebuild the paragraph: four proofErr markers followed by the same run\n# content as before (name / \" \" / \"= \" / \"World\" / \",\") but without the\n# leftover pPr + empty run.\n$replacementXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' `\n    + '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' `\n    + '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' `\n    + '<pkg:xmlData>' `\n    + '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' `\n    + '<w:body>' `\n    + $pOpenTag `\n    + '<w:proofErr w:type=\"spellStart\"/>' `\n    + '<w:proofErr w:type=\"spellEnd\"/>' `\n    + '<w:proofErr w:type=\"gramStart\"/>' `\n    + '<w:proofErr w:type=\"gramEnd\"/>' `\n    + '<w:r><w:t>name</w:t></w:r>' `\n    + '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' `\n    + '<w:r><w:t xml:space=\"preserve\">= </w:t></w:r>' `\n    + '<w:r><w:t>World</w:t></w:r>' `\n    + '<w:r><w:t>,</w:t></w:r>' `\n    + '</w:p>' `\n    + '</w:body>' `\n    + '</w:document>' `\n    + '</pkg:xmlData>' `\n    + '</pkg:part>' `\n    + '</pkg:package>'\n\n$target.InsertXML($replacementXml) | Out-Null\n"}
